$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 162, pushing the existing rows 162:173 down
# to 163:174 (and carrying their formatting, i.e. the date style on
# column D, along with them).
$ws.Rows("162:162").Insert()

# Populate the newly inserted row 162 with the latest weekly price
# observation for "Pepino ensalada" at Vega Monumental Concepción.
$ws.Cells.Item(162, 1).Value = 11
$ws.Cells.Item(162, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(162, 3).Value = "Bíobío"
$ws.Cells.Item(162, 4).Value = 44931
$ws.Cells.Item(162, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(162, 5).Value = 8
$ws.Cells.Item(162, 6).Value = 100112043
$ws.Cells.Item(162, 7).Value = "Pepino ensalada"
$ws.Cells.Item(162, 8).Value = "Sin especificar"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 180
$ws.Cells.Item(162, 11).Value = 14000
$ws.Cells.Item(162, 12).Value = 15000
$ws.Cells.Item(162, 13).Value = 14444
$ws.Cells.Item(162, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(162, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(162, 16).Value = 241
$ws.Cells.Item(162, 17).Value = 60
$ws.Cells.Item(162, 18).Value = "Hortaliza"
